# The data rows 5-13 (9 rows) of this "Sandia" subconjunto get cyclically
# rotated: the content that used to live in rows 9-13 moves up to rows 5-9,
# and the content that used to live in rows 5-8 moves down to rows 10-13.
# Only columns D (Fecha), H (Variedad), I (Calidad), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) actually change value; everything else is identical on
# every row already, so it doesn't matter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 5
$lastRow = 13
$rowCount = $lastRow - $firstRow + 1
$shift = 4  # row r takes the old values of row r+4 (wrapping within the block)

$cols = @("D", "H", "I", "J", "K", "L", "M", "P")

# Snapshot the current ("before") values for the columns that change.
# Use Value2() (called as a method, not a bare property) so that dates come
# back as their raw serial number instead of a formatted string.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2()
    }
    $snapshot[$r] = $rowVals
}

# Write back the rotated values.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $r + $shift
    if ($srcRow -gt $lastRow) {
        $srcRow = $srcRow - $rowCount
    }
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $srcVals[$col]
    }
}
